# Insert a new weekly data row for "Ají" (Hortaliza) right after the
# existing row 1159, pushing all subsequent rows down by one.
# This corresponds to the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 1160 (existing rows 1160..1243 shift to 1161..1244)
$ws.Rows.Item(1160).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(1160, 1).Value  = 6
$ws.Cells.Item(1160, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1160, 3).Value  = "Metropolitana"
$ws.Cells.Item(1160, 4).Value  = 45021
$ws.Cells.Item(1160, 5).Value  = 13
$ws.Cells.Item(1160, 6).Value  = 100112021
$ws.Cells.Item(1160, 7).Value  = "Ají"
$ws.Cells.Item(1160, 8).Value  = "Americana (o)"
$ws.Cells.Item(1160, 9).Value  = "Primera"
$ws.Cells.Item(1160, 10).Value = 580
$ws.Cells.Item(1160, 11).Value = 17000
$ws.Cells.Item(1160, 12).Value = 18000
$ws.Cells.Item(1160, 13).Value = 17552
$ws.Cells.Item(1160, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(1160, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1160, 16).Value = 702
$ws.Cells.Item(1160, 17).Value = 25
$ws.Cells.Item(1160, 18).Value = "Hortaliza"
